$d = $word.ActiveDocument
$s = $d.Styles("Title")
$s.XML = "<test/>"
Write-Host "done"
